# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1, right after the existing "Unnamed: 28" (AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the rest of the header row (bold font, thin box border,
# centered / top-aligned text) since Excel's COM model has no direct
# "copy style from another range" property in this bridge.
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous (thin box border)

# Every player on this roster shares the team's 2007 season record:
# 82 wins, 80 losses, 0 ties.
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 82  # AD
    $ws.Cells.Item($r, 31).Value = 80  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
